$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.575.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.30%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.612.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.26%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'531.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.87%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'142.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.70%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = "'0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.10%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +7.11%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.37%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.334"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.61%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.32%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.082.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.72%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'58.523.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.28%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'20.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.76%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.624.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.71%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.72%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.63%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'334.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.07%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'66.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.86%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.412"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.44%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.162"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.56%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.50%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'USDe"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.04%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'PEPE"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0731"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.23%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -3.13%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'5.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.46%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'151.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.89%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'18.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.74%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.88%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.92%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.826"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.31%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.817"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.40%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.98%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.29%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'280.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.66%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.11%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'WhiteBITCoin"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'10.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Mantle"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.591"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.68%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0531"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0936"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.07%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'18.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.88%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0223"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.41%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.937.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.79%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'4.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.93%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'17.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.46%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'111.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.54%  "
$ws.Range("E51").Style = "Normal"
